$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the embedded version string in A1
$ws.Range("A1").Value = "row_definition=8&row_data_start=9&skip_empty_col=20&version=2.0.3"

# Update the "Exported:" timestamp in A4
$ws.Range("A4").Value = "Exported: 2016-07-23 02:09:43"

# Insert two new data rows before the current "EndOfData" marker row (row 13),
# which pushes it down to row 15.
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()

# Clone formatting of the existing data row (row 10) into the two new rows,
# then overwrite the cell values with the new sitemap test entries.
$ws.Range("A10:S10").Copy($ws.Range("A11:S11"))
$ws.Range("A10:S10").Copy($ws.Range("A12:S12"))

# Row 11: "Path TEST 0 bite"
$ws.Cells.Item(11, 3).Value = "Path TEST 0 bite"
$ws.Cells.Item(11, 11).Value = "/path_test_0bite/"
$ws.Cells.Item(11, 12).Value = "/0bite.html"

# Row 12: "Path TEST not exists"
$ws.Cells.Item(12, 3).Value = "Path TEST not exists"
$ws.Cells.Item(12, 11).Value = "/path_test_not_exists/"
$ws.Cells.Item(12, 12).Value = "/not_exists.html"
